# Apply the "5th commit - corrected the otp auth" data update:
# Append 4 new donor rows to the "Donors" sheet (rows 12-15) and
# 1 new patient row to the "Patients" sheet (row 8).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Donors sheet
# ---------------------------------------------------------------
$donors = $wb.Worksheets.Item("Donors")

$donorRows = @(
    @("690ee29e6dfffb265cf88022", "Ashwini Shenoy B", "ashwinishenoyb@gmail.com", "7026438371", "A+", "My Current Location", 77.6208384, 12.9564672, $false, $false, $false, 45969.497787939814),
    @("690ee498bdd4773c2d8f92a9", "Ashwini Shenoy B", "shenoybashwini@gmail.com", "7026438371", "O+", "My Current Location", 77.6208384, 12.9564672, $false, $false, $true,  45969.50364131945),
    @("690ee579bdd4773c2d8f92ae", "Ashwini Shenoy B", "shenoybashwini@gmail.com", "7026438371", "O+", "My Current Location", 77.6208384, 12.9564672, $false, $false, $true,  45969.50624263889),
    @("690ee87057df1e2407c8ae77", "Ashwini Shenoy B", "ashenoyb@gmail.com",       "1234567890", "A+", "My Current Location", 75.1239547, 15.3647083, $false, $false, $true,  45969.51502623843)
)

$startRow = 12
for ($i = 0; $i -lt $donorRows.Count; $i++) {
    $r = $startRow + $i
    $row = $donorRows[$i]

    $donors.Cells.Item($r, 1).Value = $row[0]
    $donors.Cells.Item($r, 2).Value = $row[1]
    $donors.Cells.Item($r, 3).Value = $row[2]

    # Column D (phone) is a digit-only string; force it to stay text
    # (otherwise it auto-converts to a number) the same way the existing
    # D10/D11 phone numbers are stored.
    $donors.Cells.Item($r, 4).NumberFormat = "@"
    $donors.Cells.Item($r, 4).Value = $row[3]
    $donors.Cells.Item($r, 4).ClearFormats()

    $donors.Cells.Item($r, 5).Value = $row[4]
    $donors.Cells.Item($r, 6).Value = $row[5]
    $donors.Cells.Item($r, 7).Value = $row[6]
    $donors.Cells.Item($r, 8).Value = $row[7]
    $donors.Cells.Item($r, 9).Value = $row[8]
    $donors.Cells.Item($r, 10).Value = $row[9]
    $donors.Cells.Item($r, 11).Value = $row[10]
    $donors.Cells.Item($r, 12).Value = $row[11]
}

# The timestamp in the last donor row keeps the date-number-format style
# that used to sit on L11 (Excel drags the format down when a row is
# added right after it) - move (not copy) the style from L11 to L15.
$donors.Cells.Item(11, 12).Copy()
$donors.Cells.Item(15, 12).PasteSpecial(-4122)  # xlPasteFormats
$donors.Cells.Item(11, 12).ClearFormats()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Patients sheet
# ---------------------------------------------------------------
$patients = $wb.Worksheets.Item("Patients")

$patients.Cells.Item(8, 1).Value  = "690ee6232f7207491479c1a9"
$patients.Cells.Item(8, 2).Value  = "Ashwini Shenoy B"
$patients.Cells.Item(8, 3).Value  = "ashenoyb@gmail.com"

$patients.Cells.Item(8, 4).NumberFormat = "@"
$patients.Cells.Item(8, 4).Value = "7026438371"
$patients.Cells.Item(8, 4).ClearFormats()

$patients.Cells.Item(8, 5).Value  = "A+"
$patients.Cells.Item(8, 6).Value  = "mangalore"
$patients.Cells.Item(8, 7).Value  = 1
$patients.Cells.Item(8, 8).Value  = 1
$patients.Cells.Item(8, 9).Value  = "Medium"
$patients.Cells.Item(8, 10).Value = 45969.5082124537
